$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values (shared strings): "Full name" -> "full_name", "Email" -> "email"
$ws.Range("A1").Value = "full_name"
$ws.Range("B1").Value = "email"

# Update the active selection on the sheet to B5
$ws.Range("B5").Select()
